# Update "EC" (Estado de Cuenta) data for NIT 9006687626: replaces the 37
# "Periodo Mora" detail rows (16-52) with a freshly ordered period list
# (ascending 1703 -> 2003 instead of the previous descending 2003 -> 1703)
# and refreshes the "Valor Mora" / accrued-value figures for each row, per
# the "Actualiza base de datos EC y agrega parte 1" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# row -> (Periodo Mora text, Valor Mora (F), Salario Basico-derived total (G))
$rowData = @{
    16 = @{ Periodo = "1703"; F = 27578; G = 781242 }
    17 = @{ Periodo = "1704"; F = 27578; G = 781242 }
    18 = @{ Periodo = "1705"; F = 27578; G = 781242 }
    19 = @{ Periodo = "1706"; F = 27578; G = 781242 }
    20 = @{ Periodo = "1707"; F = 27578; G = 781242 }
    21 = @{ Periodo = "1708"; F = 27578; G = 781242 }
    22 = @{ Periodo = "1709"; F = 27578; G = 781242 }
    23 = @{ Periodo = "1710"; F = 27578; G = 781242 }
    24 = @{ Periodo = "1711"; F = 27578; G = 781242 }
    25 = @{ Periodo = "1712"; F = 27578; G = 781242 }
    26 = @{ Periodo = "1801"; F = 27578; G = 781242 }
    27 = @{ Periodo = "1802"; F = 27578; G = 781242 }
    28 = @{ Periodo = "1803"; F = 27578; G = 781242 }
    29 = @{ Periodo = "1804"; F = 27578; G = 781242 }
    30 = @{ Periodo = "1805"; F = 27578; G = 781242 }
    31 = @{ Periodo = "1806"; F = 27578; G = 781242 }
    32 = @{ Periodo = "1807"; F = 27578; G = 781242 }
    33 = @{ Periodo = "1808"; F = 27578; G = 781242 }
    34 = @{ Periodo = "1809"; F = 27578; G = 781242 }
    35 = @{ Periodo = "1810"; F = 31249; G = 781242 }
    36 = @{ Periodo = "1811"; F = 31249; G = 781242 }
    37 = @{ Periodo = "1812"; F = 31249; G = 781242 }
    38 = @{ Periodo = "1901"; F = 31249; G = 781242 }
    39 = @{ Periodo = "1902"; F = 31249; G = 781242 }
    40 = @{ Periodo = "1903"; F = 31249; G = 781242 }
    41 = @{ Periodo = "1904"; F = 31249; G = 781242 }
    42 = @{ Periodo = "1905"; F = 31249; G = 781242 }
    43 = @{ Periodo = "1906"; F = 31249; G = 781242 }
    44 = @{ Periodo = "1907"; F = 31249; G = 781242 }
    45 = @{ Periodo = "1908"; F = 31249; G = 781242 }
    46 = @{ Periodo = "1909"; F = 31249; G = 781242 }
    47 = @{ Periodo = "1910"; F = 31249; G = 781242 }
    48 = @{ Periodo = "1911"; F = 31249; G = 781242 }
    49 = @{ Periodo = "1912"; F = 31249; G = 781242 }
    50 = @{ Periodo = "2001"; F = 31249; G = 781242 }
    51 = @{ Periodo = "2002"; F = 31249; G = 781242 }
    52 = @{ Periodo = "2003"; F = 31249; G = 781242 }
}

foreach ($row in $rowData.Keys) {
    $data = $rowData[$row]
    $ws.Cells.Item($row, 5).Value = $data.Periodo   # E: Periodo Mora
    $ws.Cells.Item($row, 6).Value = $data.F         # F: Valor Mora
    $ws.Cells.Item($row, 7).Value = $data.G         # G: accrued total
}

